$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set header row values (row 1) for new columns E:G
$ws.Range("E1").Value = "Exceeded 1MB"
$ws.Range("F1").Value = "Incorrect Dimensions"
$ws.Range("G1").Value = "Unsupported File Format"

# Set data row values (row 2) for new columns E:G
$ws.Range("E2").Value = "D:\MARINA\SDET\GroupProject\Exceeded 1MB.jpg"
$ws.Range("F2").Value = "D:\MARINA\SDET\GroupProject\Incorrect Dimensions.jpg"
$ws.Range("G2").Value = "D:\MARINA\SDET\GroupProject\Unsupported FIle Format.bmp"

# Copy formatting from the existing header/data cells onto the new ones
# so the new cells reuse the same cell styles (bold/border header, centered/border data)
$ws.Range("D1").Copy()
$ws.Range("E1:G1").PasteSpecial(-4122)

$ws.Range("D2").Copy()
$ws.Range("E2:G2").PasteSpecial(-4122)

# Adjust column widths: A:C narrower (13), D:G wider (~52.9)
$ws.Range("A1:C2").ColumnWidth = 12.166666666666666
$ws.Range("D1:G2").ColumnWidth = 52

# Update selection to D10
$ws.Range("D10").Select()
